$wb = $excel.ActiveWorkbook

# --- Sheet "2025" (row 2 values updated from server) ---
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("A2").Value = 0.1996117116393296
$ws2025.Range("B2").Value = 0.09147347743582362
$ws2025.Range("E2").Value = 0.196337161233252
$ws2025.Range("I2").Value = 0.9093716833333334
$ws2025.Range("M2").Value = 0.08582791666666666
$ws2025.Range("N2").Value = 8.076984323356477
$ws2025.Range("O2").Value = 5.602091941527175

# --- Sheet "2030" (row 2 values updated from server) ---
$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("B2").Value = 0.01569748571225271
$ws2030.Range("E2").Value = 0.3501989933240413
$ws2030.Range("I2").Value = 0.7880958234314523
$ws2030.Range("M2").Value = 0.08384033333333336
$ws2030.Range("N2").Value = 13.54405293432583
$ws2030.Range("O2").Value = 4.112137734330497

# --- Sheet "2035" (row 2 values updated from server) ---
$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("E2").Value = 0.2898581252785336
$ws2035.Range("G2").Value = 0.1190881252125949
$ws2035.Range("I2").Value = 0.4976846007978426
$ws2035.Range("L2").Value = 0.06366467577070564
$ws2035.Range("M2").Value = 0.04822266666666594
$ws2035.Range("N2").Value = 4.564495481057364
$ws2035.Range("O2").Value = 2.943764624981792
